# Change the exercise bullet from "Install samtools on a VM" to
# "Install docker and python3.7 on a VM" (slide 3, "Content Placeholder 2").
#
# The paragraph is made up of three runs:
#   "Install " + "samtools" + " on a VM"
# Each piece is rewritten in place via Characters()-based substring
# replacement so the surrounding runs (and their rPr formatting) stay
# untouched, matching the original run layout:
#   "Install docker " + "and python3.7 " + "on a VM"

$p = $ppt.ActivePresentation

# Locate the shape that holds the "samtools" bullet instead of assuming a
# fixed slide/shape index.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text.IndexOf("samtools") -ge 0) {
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Helper: find the paragraph (1-based index) whose text contains $needle.
function Find-ParagraphIndex($textRange, $needle) {
    $count = $textRange.Paragraphs().Count
    for ($pi = 1; $pi -le $count; $pi++) {
        if ($textRange.Paragraphs($pi, 1).Text.IndexOf($needle) -ge 0) {
            return $pi
        }
    }
    return -1
}

$paraIdx = Find-ParagraphIndex $tr "samtools"

# --- Replace the middle run: "samtools" -> "and python3.7 " -------------
$para = $tr.Paragraphs($paraIdx, 1)
$oldMid = "samtools"
$midStart = $para.Text.IndexOf($oldMid) + 1
$runMid = $para.Characters($midStart, $oldMid.Length)
$runMid.Text = "and python3.7 "

# --- Replace the leading run: "Install " -> "Install docker " ----------
$para = $tr.Paragraphs($paraIdx, 1)
$prefix = "Install "
$runPre = $para.Characters(1, $prefix.Length)
$runPre.Text = "Install docker "

# --- Replace the trailing run: " on a VM" -> "on a VM" -----------------
$para = $tr.Paragraphs($paraIdx, 1)
$suffix = " on a VM"
$sufStart = $para.Text.LastIndexOf($suffix) + 1
$runSuf = $para.Characters($sufStart, $suffix.Length)
$runSuf.Text = "on a VM"

$para = $tr.Paragraphs($paraIdx, 1)
Write-Output "Updated bullet text: $($para.Text)"
